# Auto-generated PowerShell COM-interop script
# Applies the cryptos.xlsx price/volume update described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $range = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "1.012") are not
    # auto-converted to numbers, matching the original inline-string cell type.
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    # Restore the default style so we do not leave a stray text-format style
    # index behind (the source cells carried no explicit style).
    $range.Style = "Normal"
}

Set-TextCellValue "D2" "27.458.77"
Set-TextCellValue "E2" "  +2.04%  "
Set-TextCellValue "D3" "1.837.86"
Set-TextCellValue "E3" "  +1.26%  "
Set-TextCellValue "D4" "1.012"
Set-TextCellValue "E4" "  +1.06%  "
Set-TextCellValue "D5" "314.38"
Set-TextCellValue "E5" "  +1.59%  "
Set-TextCellValue "E6" "  +0.95%  "
Set-TextCellValue "D7" "0.4736"
Set-TextCellValue "E7" "  +1.73%  "
Set-TextCellValue "D8" "0.3693"
Set-TextCellValue "E8" "  +0.76%  "
Set-TextCellValue "D9" "0.07460"
Set-TextCellValue "E9" "  +1.45%  "
Set-TextCellValue "D10" "0.8854"
Set-TextCellValue "E10" "  +1.91%  "
Set-TextCellValue "D11" "20.44"
Set-TextCellValue "E11" "  +0.60%  "
Set-TextCellValue "D12" "1.919.52"
Set-TextCellValue "E12" "  +3.34%  "
Set-TextCellValue "D13" "0.07329"
Set-TextCellValue "E13" "  +3.49%  "
Set-TextCellValue "D14" "5.453"
Set-TextCellValue "E14" "  +1.28%  "
Set-TextCellValue "E15" "  +1.80%  "
Set-TextCellValue "E16" "  +1.12%  "
Set-TextCellValue "D17" "1.008"
Set-TextCellValue "E17" "  +0.54%  "
Set-TextCellValue "D18" "0.000008822"
Set-TextCellValue "E18" "  +1.32%  "
Set-TextCellValue "E20" "  +1.02%  "
Set-TextCellValue "D21" "27.487.38"
Set-TextCellValue "E21" "  +2.01%  "
Set-TextCellValue "D22" "5.327"
Set-TextCellValue "E22" "  +0.51%  "
Set-TextCellValue "D23" "10.70"
Set-TextCellValue "E23" "  +0.52%  "
Set-TextCellValue "D24" "2.130.85"
Set-TextCellValue "E24" "  +2.72%  "
Set-TextCellValue "E25" "  +0.78%  "
Set-TextCellValue "E26" "  +0.79%  "
Set-TextCellValue "D27" "18.61"
Set-TextCellValue "E27" "  +1.34%  "
Set-TextCellValue "D28" "2.150"
Set-TextCellValue "E28" "  +0.37%  "
Set-TextCellValue "D29" "5.252"
Set-TextCellValue "E29" "  -0.07%  "
Set-TextCellValue "D30" "117.90"
Set-TextCellValue "E30" "  +2.18%  "
Set-TextCellValue "D31" "0.08997"
Set-TextCellValue "E31" "  +0.86%  "
Set-TextCellValue "D32" "0.7573"
Set-TextCellValue "E32" "  +0.15%  "
Set-TextCellValue "D33" "1.181"
Set-TextCellValue "E33" "  +2.20%  "
Set-TextCellValue "D34" "4.557"
Set-TextCellValue "E34" "  +1.61%  "
Set-TextCellValue "D35" "2.937"
Set-TextCellValue "E35" "  +0.89%  "
Set-TextCellValue "D36" "1.013"
Set-TextCellValue "E36" "  +1.13%  "
Set-TextCellValue "D37" "1.106"
Set-TextCellValue "E37" "  +1.95%  "
Set-TextCellValue "D38" "0.05337"
Set-TextCellValue "E38" "  +1.17%  "
Set-TextCellValue "E39" "  +0.45%  "
Set-TextCellValue "D40" "2.993"
Set-TextCellValue "D41" "7.326"
Set-TextCellValue "E41" "  +1.08%  "
Set-TextCellValue "D42" "2.410"
Set-TextCellValue "E42" "  +5.43%  "
Set-TextCellValue "D43" "0.5334"
Set-TextCellValue "E43" "  +0.55%  "
Set-TextCellValue "D44" "0.1661"
Set-TextCellValue "E44" "  +0.49%  "
Set-TextCellValue "D45" "8.516"
Set-TextCellValue "E45" "  +0.99%  "
Set-TextCellValue "D46" "0.4912"
Set-TextCellValue "E46" "  +0.76%  "
Set-TextCellValue "D47" "10.52"
Set-TextCellValue "E47" "  +1.31%  "
Set-TextCellValue "E48" "  +1.11%  "
Set-TextCellValue "D49" "105.12"
Set-TextCellValue "E49" "  +1.80%  "
Set-TextCellValue "E50" "  +1.14%  "
Set-TextCellValue "D51" "0.06311"
Set-TextCellValue "E51" "  +0.30%  "
